$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "12"
$ws.Range("J6").Value = "先遣侦查，保持距离发射腐质胆汁。"

$ws.Range("E7").Value = "30"
$ws.Range("G7").Value = "10"
$ws.Range("J7").Value = "歌祭徒吟唱增幅附近同伴。"

$ws.Range("D8").Value = "75"
$ws.Range("G8").Value = "8"
$ws.Range("J8").Value = "雾潮袭来后出现呼嚎者，注意理智流失。"

$ws.Range("D9").Value = "115"
$ws.Range("E9").Value = "35"
$ws.Range("G9").Value = "6"
$ws.Range("J9").Value = "虚壳哨兵扫描横扫，伴随餍爬者挤压。"

$ws.Range("D10").Value = "160"
$ws.Range("G10").Value = "14"
$ws.Range("J10").Value = "碎影成群突进，逼迫频繁拉扯走位。"

$ws.Range("D11").Value = "210"
$ws.Range("E11").Value = "40"
$ws.Range("G11").Value = "8"
$ws.Range("J11").Value = "掘锚者直冲核心，需迅速打断。"
